# Updated cryptos list on Wed Apr 26 14:51:11 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (D) and "Volume(1h)" (E) columns for every coin row
# with a new snapshot pulled from coinranking.com, plus re-rank a couple of
# coins that swapped positions (RenderToken/Decentraland and
# Cronos/NEARProtocol), which carries their Coin name (B) and Link (C)
# along with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal TEXT (matches the workbook's
# original inline-string cells) even when it looks like a number
# (e.g. "1.001", "12.60", "0.06670") so Excel doesn't silently reinterpret
# it as a numeric value and drop significant trailing/format digits.
# NumberFormat is restored to the default afterwards so the cell's style
# is left exactly as it was found.
function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
$ws.Cells.Item(2, 4).Value = "29.693.34"
$ws.Cells.Item(2, 5).Value = "  +8.49%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.944.99"
$ws.Cells.Item(3, 5).Value = "  +7.07%  "

# Row 4
Set-TextValue 4 4 "1.001"
$ws.Cells.Item(4, 5).Value = "  -0.26%  "

# Row 5
Set-TextValue 5 4 "341.56"
$ws.Cells.Item(5, 5).Value = "  +3.22%  "

# Row 6
Set-TextValue 6 4 "1.001"
$ws.Cells.Item(6, 5).Value = "  -0.16%  "

# Row 7
Set-TextValue 7 4 "0.4772"
$ws.Cells.Item(7, 5).Value = "  +4.59%  "

# Row 8
Set-TextValue 8 4 "0.4133"
$ws.Cells.Item(8, 5).Value = "  +8.58%  "

# Row 9
Set-TextValue 9 4 "48.42"
$ws.Cells.Item(9, 5).Value = "  +5.46%  "

# Row 10
Set-TextValue 10 4 "0.08245"
$ws.Cells.Item(10, 5).Value = "  +5.17%  "

# Row 11
Set-TextValue 11 4 "1.042"
$ws.Cells.Item(11, 5).Value = "  +8.80%  "

# Row 12
Set-TextValue 12 4 "22.64"

# Row 13
$ws.Cells.Item(13, 4).Value = "1.940.44"
$ws.Cells.Item(13, 5).Value = "  +5.53%  "

# Row 14
Set-TextValue 14 4 "6.182"
$ws.Cells.Item(14, 5).Value = "  +5.90%  "

# Row 15
Set-TextValue 15 4 "7.421"
$ws.Cells.Item(15, 5).Value = "  +5.21%  "

# Row 16
Set-TextValue 16 4 "92.42"

# Row 17
Set-TextValue 17 4 "1.002"
$ws.Cells.Item(17, 5).Value = "  -0.15%  "

# Row 18
Set-TextValue 18 4 "0.00001064"
$ws.Cells.Item(18, 5).Value = "  +4.51%  "

# Row 19
Set-TextValue 19 4 "0.06670"
$ws.Cells.Item(19, 5).Value = "  +1.47%  "

# Row 20
Set-TextValue 20 4 "18.07"
$ws.Cells.Item(20, 5).Value = "  +5.82%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -0.16%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "29.658.11"
$ws.Cells.Item(22, 5).Value = "  +8.40%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +6.33%  "

# Row 24
Set-TextValue 24 4 "11.24"
$ws.Cells.Item(24, 5).Value = "  +4.21%  "

# Row 25
Set-TextValue 25 4 "2.282"

# Row 26
$ws.Cells.Item(26, 4).Value = "2.164.92"
$ws.Cells.Item(26, 5).Value = "  +5.42%  "

# Row 27
Set-TextValue 27 4 "160.52"
$ws.Cells.Item(27, 5).Value = "  +3.20%  "

# Row 28
Set-TextValue 28 4 "20.18"
$ws.Cells.Item(28, 5).Value = "  +4.90%  "

# Row 29
Set-TextValue 29 4 "2.198"
$ws.Cells.Item(29, 5).Value = "  +7.91%  "

# Row 30
Set-TextValue 30 4 "5.636"
$ws.Cells.Item(30, 5).Value = "  +7.59%  "

# Row 31
Set-TextValue 31 4 "122.34"
$ws.Cells.Item(31, 5).Value = "  +4.16%  "

# Row 32
Set-TextValue 32 4 "1.027"
$ws.Cells.Item(32, 5).Value = "  +10.28%  "

# Row 33
Set-TextValue 33 4 "0.09652"
$ws.Cells.Item(33, 5).Value = "  +4.13%  "

# Row 34
Set-TextValue 34 4 "1.471"
$ws.Cells.Item(34, 5).Value = "  +12.22%  "

# Row 35
Set-TextValue 35 4 "3.682"
$ws.Cells.Item(35, 5).Value = "  +3.35%  "

# Row 36
Set-TextValue 36 4 "5.493"
$ws.Cells.Item(36, 5).Value = "  +5.54%  "

# Row 37
Set-TextValue 37 4 "0.06294"
$ws.Cells.Item(37, 5).Value = "  +6.67%  "

# Row 38
Set-TextValue 38 4 "0.02328"
$ws.Cells.Item(38, 5).Value = "  +6.87%  "

# Row 39
Set-TextValue 39 4 "8.602"
$ws.Cells.Item(39, 5).Value = "  +6.55%  "

# Row 40
Set-TextValue 40 4 "1.194"
$ws.Cells.Item(40, 5).Value = "  +5.16%  "

# Row 41
Set-TextValue 41 4 "0.6107"
$ws.Cells.Item(41, 5).Value = "  +6.60%  "

# Row 42
Set-TextValue 42 4 "10.71"

# Row 43
Set-TextValue 43 4 "0.1905"
$ws.Cells.Item(43, 5).Value = "  +5.15%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -0.11%  "

# Row 45
Set-TextValue 45 4 "1.271"
$ws.Cells.Item(45, 5).Value = "  +0.42%  "

# Row 46
Set-TextValue 46 4 "12.60"
$ws.Cells.Item(46, 5).Value = "  +6.62%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "RenderToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 47 4 "2.362"
$ws.Cells.Item(47, 5).Value = "  +32.51%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "Decentraland"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue 48 4 "0.5717"
$ws.Cells.Item(48, 5).Value = "  +6.38%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue 49 4 "0.07437"
$ws.Cells.Item(49, 5).Value = "  +13.19%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "NEARProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue 50 4 "1.999"
$ws.Cells.Item(50, 5).Value = "  +7.37%  "

# Row 51
Set-TextValue 51 4 "114.41"
$ws.Cells.Item(51, 5).Value = "  +4.58%  "
